$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: insert "Cohort" into A1, shift existing headers right,
# and add new "Participants" header in F1.
$ws.Range("A1").Value = "Cohort"
$ws.Range("B1").Value = "Component"
$ws.Range("C1").Value = "Direction"
$ws.Range("D1").Value = "Mean"
$ws.Range("E1").Value = "Standard_Deviation"
$ws.Range("F1").Value = "Participants"

# Copy header style (A1) across the new header cells
$ws.Range("A1").Copy()
$ws.Range("B1:F1").PasteSpecial(-4122)  # xlPasteFormats

$data = @(
    @("100s", "Rambling",  "X", -0.03658395628494156,              0.3931959601114328,              18),
    @("100s", "Rambling",  "Y", -0.08625170631061047,              1.945771240426168,               18),
    @("100s", "Trembling", "X", [double]"-8.454678618618019E-07",  [double]"4.245181673844462E-05",  18),
    @("100s", "Trembling", "Y", [double]"2.000999336161684E-05",   0.000452942851766289,             18),
    @("200s", "Rambling",  "X", 0.09327416293317954,               0.4596924592754691,               20),
    @("200s", "Rambling",  "Y", 0.1039070088240599,                1.850877978605084,                20),
    @("200s", "Trembling", "X", [double]"-7.822513613717468E-07",  [double]"4.342558459508203E-05",  20),
    @("200s", "Trembling", "Y", [double]"-3.048420284778794E-05",  0.0007241024510988817,            20)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $row++
}
